# chore: update Sheets via scheduled runner
# Refreshes currentAveragePrice / LevePrice / LeveProfit columns (H:N) for a
# handful of leve rows across several crafting-job sheets, based on
# refreshed Universalis-style market data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2606272.5
$ws.Range("I40").Value = 6251350
$ws.Range("J40").Value = 2645.7144
$ws.Range("K40").Value = 6251350
$ws.Range("L40").Value = 2645.7144
$ws.Range("M40").Value = -6251175
$ws.Range("N40").Value = -2995.7144
$ws.Range("H116").Value = 2313.6365
$ws.Range("I116").Value = 1762.5
$ws.Range("J116").Value = 3783.3333
$ws.Range("K116").Value = 1762.5
$ws.Range("L116").Value = 3783.3333
$ws.Range("M116").Value = 1679.5
$ws.Range("N116").Value = -10667.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H57").Value = 3000
$ws.Range("I57").Value = 3000
$ws.Range("K57").Value = 3000
$ws.Range("M57").Value = -2516

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 35498.11
$ws.Range("J109").Value = 35498.11
$ws.Range("L109").Value = 35498.11
$ws.Range("N109").Value = -38272.11

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 251.375
$ws.Range("I22").Value = 242.2
$ws.Range("J22").Value = 266.66666
$ws.Range("K22").Value = 242.2
$ws.Range("L22").Value = 266.66666
$ws.Range("M22").Value = 107.8
$ws.Range("N22").Value = -966.66666
$ws.Range("H31").Value = 2678.0508
$ws.Range("I31").Value = 2297.4583
$ws.Range("J31").Value = 2939.0286
$ws.Range("K31").Value = 2297.4583
$ws.Range("L31").Value = 2939.0286
$ws.Range("M31").Value = -2002.4583
$ws.Range("N31").Value = -3529.0286
$ws.Range("H34").Value = 2678.0508
$ws.Range("I34").Value = 2297.4583
$ws.Range("J34").Value = 2939.0286
$ws.Range("K34").Value = 2297.4583
$ws.Range("L34").Value = 2939.0286
$ws.Range("M34").Value = -2095.4583
$ws.Range("N34").Value = -3343.0286
$ws.Range("H75").Value = 20238.715
$ws.Range("J75").Value = 20238.715
$ws.Range("L75").Value = 20238.715
$ws.Range("N75").Value = -22234.715
$ws.Range("H78").Value = 20238.715
$ws.Range("J78").Value = 20238.715
$ws.Range("L78").Value = 60716.145
$ws.Range("N78").Value = -70700.145
$ws.Range("H98").Value = 21000
$ws.Range("I98").Value = 10000
$ws.Range("K98").Value = 10000
$ws.Range("M98").Value = -7754
$ws.Range("H127").Value = 10000
$ws.Range("J127").Value = 10000
$ws.Range("L127").Value = 10000
$ws.Range("N127").Value = -19920
$ws.Range("H141").Value = 56289.766
$ws.Range("I141").Value = 16666.666
$ws.Range("J141").Value = 64780.43
$ws.Range("K141").Value = 16666.666
$ws.Range("L141").Value = 64780.43
$ws.Range("M141").Value = -11486.666
$ws.Range("N141").Value = -75140.42999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 416.66666
$ws.Range("H113").Value = 6790665
$ws.Range("I113").Value = 6944901.5
$ws.Range("J113").Value = 6667276
$ws.Range("K113").Value = 20834704.5
$ws.Range("L113").Value = 20001828
$ws.Range("M113").Value = -20832534.5
$ws.Range("N113").Value = -20006168
$ws.Range("H129").Value = 15873949
$ws.Range("I129").Value = 720.8333
$ws.Range("J129").Value = 37038252
$ws.Range("K129").Value = 2162.4999
$ws.Range("L129").Value = 111114756
$ws.Range("M129").Value = 2837.5001
$ws.Range("N129").Value = -111124756
$ws.Range("H131").Value = 799.86
$ws.Range("I131").Value = 415
$ws.Range("J131").Value = 842.6222
$ws.Range("K131").Value = 1245
$ws.Range("L131").Value = 2527.8666
$ws.Range("M131").Value = 3795
$ws.Range("N131").Value = -12607.8666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 5000
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 5000
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 5000
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -6166
$ws.Range("H69").Value = 17999.5
$ws.Range("J69").Value = 17999.5
$ws.Range("L69").Value = 17999.5
$ws.Range("N69").Value = -19497.5
$ws.Range("H72").Value = 17999.5
$ws.Range("J72").Value = 17999.5
$ws.Range("L72").Value = 53998.5
$ws.Range("N72").Value = -61486.5
$ws.Range("H126").Value = 3755.4167
$ws.Range("I126").Value = 3028.3333
$ws.Range("K126").Value = 9084.999899999999
$ws.Range("M126").Value = -6614.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2111753.8
$ws.Range("I22").Value = 4219607.5
$ws.Range("J22").Value = 3900
$ws.Range("K22").Value = 4219607.5
$ws.Range("L22").Value = 3900
$ws.Range("M22").Value = -4219312.5
$ws.Range("N22").Value = -4490
$ws.Range("H27").Value = 2111753.8
$ws.Range("I27").Value = 4219607.5
$ws.Range("J27").Value = 3900
$ws.Range("K27").Value = 4219607.5
$ws.Range("L27").Value = 3900
$ws.Range("M27").Value = -4219500.5
$ws.Range("N27").Value = -4114
$ws.Range("H82").Value = 1600
$ws.Range("I82").Value = 1500
$ws.Range("J82").Value = 1766.6666
$ws.Range("K82").Value = 1500
$ws.Range("L82").Value = 1766.6666
$ws.Range("M82").Value = -1139
$ws.Range("N82").Value = -2488.6666
$ws.Range("H85").Value = 1600
$ws.Range("I85").Value = 1500
$ws.Range("J85").Value = 1766.6666
$ws.Range("K85").Value = 1500
$ws.Range("L85").Value = 1766.6666
$ws.Range("M85").Value = -252
$ws.Range("N85").Value = -4262.6666
$ws.Range("H93").Value = 1267.25
$ws.Range("I93").Value = 1274.3334
$ws.Range("J93").Value = 1256.625
$ws.Range("K93").Value = 1274.3334
$ws.Range("L93").Value = 1256.625
$ws.Range("M93").Value = -26.33339999999998
$ws.Range("N93").Value = -3752.625
$ws.Range("H136").Value = 4624.6
$ws.Range("I136").Value = 3484.6924
$ws.Range("J136").Value = 6184.4736
$ws.Range("K136").Value = 10454.0772
$ws.Range("L136").Value = 18553.4208
$ws.Range("M136").Value = -7904.0772
$ws.Range("N136").Value = -23653.4208

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 1575
$ws.Range("I21").Value = 150
$ws.Range("J21").Value = 3000
$ws.Range("K21").Value = 150
$ws.Range("L21").Value = 3000
$ws.Range("M21").Value = 85
$ws.Range("N21").Value = -3470
$ws.Range("H35").Value = 1575
$ws.Range("I35").Value = 150
$ws.Range("J35").Value = 3000
$ws.Range("K35").Value = 150
$ws.Range("L35").Value = 3000
$ws.Range("M35").Value = 140
$ws.Range("N35").Value = -3580
$ws.Range("H132").Value = 10883192
$ws.Range("I132").Value = 15642280
$ws.Range("J132").Value = 5276.7856
$ws.Range("K132").Value = 46926840
$ws.Range("L132").Value = 15830.3568
$ws.Range("M132").Value = -46924310
$ws.Range("N132").Value = -20890.3568
$ws.Range("H136").Value = 2882.254
$ws.Range("I136").Value = 3989.0293
$ws.Range("J136").Value = 1584.6552
$ws.Range("K136").Value = 11967.0879
$ws.Range("L136").Value = 4753.9656
$ws.Range("M136").Value = -9417.0879
$ws.Range("N136").Value = -9853.9656
